$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 22: text becomes merged string, and C/D values change
$ws.Range("B22").Value = "Pantalla de selección, Reporte OO UI Base"
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 2.5

# Update row 23: D23 gets a new value
$ws.Range("D23").Value = 2.5

# Update selection / view
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("D24").Select()
